$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The revised search-string blocks now retrieve every test-list article
# correctly, so the row that had been duplicated -- STUDY_ID 75, Layman et
# al., "Provision of ecosystem services by human-made structures in a highly
# impacted estuary" (row 80) -- is removed from the screened test list.
# Deleting the whole row shifts every following row up by one (A1:J112 ->
# A1:J111) and drops the four shared strings that were unique to that row
# (author/doi/url/title) out of sharedStrings.xml.
$ws.Rows.Item(80).Delete()

# Restore the recorded view/selection state: scrolled so row 55 is at the
# top of the window, with K78 as the active selected cell.
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("K78").Select()
